$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- Row 47: Kitchen Confidential ---
$ws.Range("A47").Value = "Kitchen Confidential"
$ws.Range("B47").Value = "Anthony Bourdain"

$ws.Range("C46").Copy()
$ws.Range("C47").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C47").Value2 = 43919
$ws.Range("D46").Copy()
$ws.Range("D47").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D47").Value2 = 43920

$ws.Range("E47").Value = "cooking;chef;restaurants;"
$ws.Range("F47").Value = "Audio"

# --- Row 48: Get Your Shit Together ---
$ws.Range("A48").Value = "Get Your Shit Together"
$ws.Range("B48").Value = "Sarah Knight"

$ws.Range("C46").Copy()
$ws.Range("C48").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C48").Value2 = 43920
$ws.Range("D46").Copy()
$ws.Range("D48").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D48").Value2 = 43920

$ws.Range("E48").Value = "self help;productivity;focus;improvement"
$ws.Range("F48").Value = "Audio"

$ws.Application.CutCopyMode = $false

# --- Lengths, filled in afterwards: row 48 first, then row 47, then fix row 46's typo ---
$ws.Range("G48").Value = "4 Hours 35 Mins"
$ws.Range("G47").Value = "8 Hours 26 Mins"
$ws.Range("G46").Value = "23 Hours 21 Mins"

# Update the view to match: scrolled down a bit, with G47 selected
$ws.Application.ActiveWindow.ScrollRow = 29
$ws.Range("G47").Select()
